$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up the vestigial row-level custom format on row 2 ---
$ws.Rows(2).ClearFormats()

# --- New shared strings / values for the shift-line summary block ---
$ws.Range("F6").Value = "Tổng thời gian"
$ws.Range("F7").Value = "Chạy"
$ws.Range("G7").Value = "Dừng"

# --- Bump base font size to 13pt across the cells that carry content ---
$ws.Range("A2").Font.Size = 13
$ws.Range("A3").Font.Size = 13
$ws.Range("A6").Font.Size = 13
$ws.Range("A7").Font.Size = 13
$ws.Range("B7").Font.Size = 13
$ws.Range("F6").Font.Size = 13

# --- "Chạy" (running) cell: green fill, white text, centered, thin border ---
$run = $ws.Range("F7")
$run.Font.Size = 13
$run.Font.ThemeColor = 2
$run.Interior.Color = 5287936
$run.HorizontalAlignment = -4108
$run.VerticalAlignment = -4108
$run.Borders.LineStyle = 1
$run.Borders.Weight = 2

# --- "Dừng" (stopped) cell: red fill, white text, centered, thin border ---
$stop = $ws.Range("G7")
$stop.Font.Size = 13
$stop.Font.ThemeColor = 2
$stop.Interior.Color = 255
$stop.HorizontalAlignment = -4108
$stop.VerticalAlignment = -4108
$stop.Borders.LineStyle = 1
$stop.Borders.Weight = 2

# --- Value cells under the Run/Stop counters: text format, centered, bordered ---
$vals = $ws.Range("F8:G8")
$vals.NumberFormat = "@"
$vals.Font.Size = 13
$vals.HorizontalAlignment = -4108
$vals.VerticalAlignment = -4108
$vals.Borders.LineStyle = 1
$vals.Borders.Weight = 2

$ws.Range("F13").Select()
